$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62, pushing the existing rows 62..154 down to 63..155.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new weekly record.
$ws.Cells.Item(62, 1).Value = 9
$ws.Cells.Item(62, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(62, 3).Value = "Metropolitana"
$ws.Cells.Item(62, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(62, 5).Value = 13
$ws.Cells.Item(62, 6).Value = 300000001
$ws.Cells.Item(62, 7).Value = "Rabanito"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 10).Value = 7900
$ws.Cells.Item(62, 11).Value = 3500
$ws.Cells.Item(62, 12).Value = 4000
$ws.Cells.Item(62, 13).Value = 3747
$ws.Cells.Item(62, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(62, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(62, 16).Value = 37
$ws.Cells.Item(62, 17).Value = 100
$ws.Cells.Item(62, 18).Value = "Hortaliza"
